# Applies the content reshuffle described by the diff: several labelled rows in
# the course-syllabus sheet (rows 10, 13-23) had their values/labels shifted and
# row 24 (the last "Requisitos" value row) was removed, collapsing the sheet from
# A1:C24 down to A1:C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos value now holds the professor name instead of the long PT objectives text ---
$ws.Range("B10").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C10").Value = "5817650 - Érica Leonor Romão"

# --- Row 13: gains an A13 label ("Programa resumido:"); B13/C13 become "Semestral"; row height 60 ---
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: label becomes "Short syllabus:"; value becomes the EN short-syllabus text ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "The environmental licensing as an instrument of environmental policy; procedures environmental licensing; licensing in the state of São Paulo: institutional aspects and applied legislation; case study during the stages of licensing."
$ws.Range("C14").Value = "The environmental licensing as an instrument of environmental policy; procedures environmental licensing; licensing in the state of São Paulo: institutional aspects and applied legislation; case study during the stages of licensing."

# --- Row 15: label becomes "Programa:"; value becomes a date placeholder; row height 120 ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: label becomes "Syllabus:"; value becomes the EN long-syllabus text ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Environmental licensing as an instrument of national environmental policy; Stages and deadlines of environmental licensing: prior license, installation and operation of the enterprise; Ventures subject to environmental licensing; Licensing in the state of São Paulo: institutional aspects, applied legislation and documentation; Technical and legal approach in environmental licensing regarding native vegetation and permanent preservation area in the State of São Paulo, case study; performance of the environmental engineer."
$ws.Range("C16").Value = "Environmental licensing as an instrument of national environmental policy; Stages and deadlines of environmental licensing: prior license, installation and operation of the enterprise; Ventures subject to environmental licensing; Licensing in the state of São Paulo: institutional aspects, applied legislation and documentation; Technical and legal approach in environmental licensing regarding native vegetation and permanent preservation area in the State of São Paulo, case study; performance of the environmental engineer."

# --- Row 17: label becomes "Avaliação:"; loses its B/C value cells entirely; row height reverts to default ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Rows.Item(17).AutoFit()

# --- Row 18: label becomes "Método:"; gains B18/C18 with the professor name; row height 60 ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("B18").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C18").Value = "5817650 - Érica Leonor Romão"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: label becomes "Critério:" (value text unchanged) ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20: label becomes "Norma de recuperação:" (value text unchanged) ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: label becomes "Bibliografia:" (value text unchanged); row height 120 ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22: label becomes "Requisitos:"; loses its B/C value cells entirely; row height reverts to default ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22:C22").Clear()
$ws.Rows.Item(22).AutoFit()

# --- Row 23: loses its A23 label; gains B23/C23 with the requisito text; row height 30 ---
$ws.Range("A23").Clear()
$ws.Range("B19").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("B23").Value = "LOB1235 -  Avaliação de Impactos Ambientais  (Requisito)`n"
$ws.Range("C23").Value = "LOB1235 -  Avaliação de Impactos Ambientais  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24 no longer exists in the target sheet; remove it ---
$ws.Rows.Item(24).Delete()

Write-Host "edit.ps1 complete"
